$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7056151032447815
$ws.Range("B1").Value = 1.34458601474762
$ws.Range("C1").Value = 4.207367897033691
$ws.Range("D1").Value = 1.531088471412659
$ws.Range("E1").Value = 0.6774513721466064
